$d = $word.ActiveDocument

# Remove the pre-existing "_GoBack" bookmark (it marks Word's last-edit
# location and will be re-created at the new edit point below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Replace "front and back end" with "proficient in all paradigms" inside
# the "Java (front and back end), Android ..." bullet.
$rng = $d.Content
[void]$rng.Find.Execute("front and back end", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "proficient in all paradigms", 2)

# Locate the freshly-inserted text so we can drop a collapsed "_GoBack"
# bookmark right after it (this also forces the surrounding text to be
# split into separate runs, matching Word's own editing behaviour).
$rng2 = $d.Content
[void]$rng2.Find.Execute("proficient in all paradigms", $true, $false, $false, $false, $false, `
                    $true, 1, $false)

$tempStart = $d.Range($rng2.Start, $rng2.Start)
$d.Bookmarks.Add("IronTempSplit", $tempStart)

$goBackPoint = $d.Range($rng2.End, $rng2.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

# The temporary bookmark was only needed to force the run split between
# "Java (" and "proficient in all paradigms"; remove it now.
$d.Bookmarks.Item("IronTempSplit").Delete()
